$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Correccion de alta de usuarios": add a new pending/fixed task right under the
# "USUARIOS" section header, flagging that user creation ("Alta de usuarios") is
# broken. This pushes the rest of that section (and everything after it) down
# by three rows, mirroring a blank-row / task-row / blank-row insertion.
$ws.Rows("52:54").Insert()

# A53: bold, 16pt, non-underlined header-ish text (distinct from the underlined
# "USUARIOS" title style) carrying the new note.
$a53 = $ws.Range("A53")
$a53.Value = "ARREGLAR ALTA DE USUARIOS "
$a53.Font.Bold = $true
$a53.Font.Size = 16
$a53.Font.Underline = $false

# B53: short explanation.
$ws.Range("B53").Value = "no funciona"

# C53: status column, marked "hecho" (done) in the same green used elsewhere
# in the Estado column.
$c53 = $ws.Range("C53")
$c53.Value = "hecho"
$c53.Font.Color = 5287936

# Reflect the author's scrolled/selected viewport after the edit.
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("C53").Select() | Out-Null
